# Adds the three new "Marzo 2019" transaction rows (94-96) to the
# "Transacciones" sheet, matching the committed data update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transacciones")

# --- Row 94 -----------------------------------------------------------
# Copy the date-format (s=1) from A93 down into the new date cells so we
# don't mint a brand-new (redundant) number-format style.
$ws.Range("A93").Copy()
$ws.Range("A94:A96").PasteSpecial(-4122)

# Copy the "savings delta" style (s=27) from P93 down the new P column.
$ws.Range("P93").Copy()
$ws.Range("P94:P96").PasteSpecial(-4122)

$ws.Range("A94").Value = 43579
$ws.Range("B94").Value = 13
$ws.Range("C94").Value = "Coca Cola"
$ws.Range("D94").Value = "Golosina"
$ws.Range("E94").Value = "Gasto"
$ws.Range("F94").Value = "Tarjeta Banamex"
$ws.Range("G94").Value = "Extra"
$ws.Range("K94").Formula = "=K93-B94"
$ws.Range("L94").Value = 2298.5700000000002
$ws.Range("M94").Value = 2
$ws.Range("N94").Formula = "=SUM(K94:M94)"
$ws.Range("O94").Formula = "=N94-4000"
$ws.Range("P94").Formula = "=O94-Ahorros!`$E`$4"

# --- Row 95 -------------------------------------------------------------
$ws.Range("A95").Value = 43580
$ws.Range("B95").Value = 200
$ws.Range("C95").Value = "Netflix"
$ws.Range("D95").Value = "Servicios"
$ws.Range("E95").Value = "Gasto"
$ws.Range("F95").Value = "Tarjeta Banamex"
$ws.Range("G95").Value = "Extra"
$ws.Range("K95").Formula = "=K94-B95"
$ws.Range("L95").Value = 2298.5700000000002
$ws.Range("M95").Value = 2
$ws.Range("N95").Formula = "=SUM(K95:M95)"
$ws.Range("O95").Formula = "=N95-4000"
$ws.Range("P95").Formula = "=O95-Ahorros!`$E`$4"

# --- Row 96 -------------------------------------------------------------
$ws.Range("A96").Value = 43581
$ws.Range("B96").Value = 24
$ws.Range("C96").Value = "Galletas "
$ws.Range("D96").Value = "Golosina"
$ws.Range("E96").Value = "Gasto"
$ws.Range("F96").Value = "Tarjeta Santander"
$ws.Range("G96").Value = "Extra"
$ws.Range("K96").Value = 6556.44
$ws.Range("L96").Formula = "=L95-B96"
$ws.Range("M96").Value = 2
$ws.Range("N96").Formula = "=SUM(K96:M96)"
$ws.Range("O96").Formula = "=N96-4000"
$ws.Range("P96").Formula = "=O96-Ahorros!`$E`$4"

# --- View state: mirror the saved selection from the source edit -------
$win = $excel.ActiveWindow
$win.ScrollRow = 84
$win.ScrollColumn = 1
$ws.Range("R96").Select()

"Added rows 94-96 to Transacciones"
